$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("input")

# Update input parameter values (C2, D3:D8)
$wsInput.Range("C2").Value = 5
$wsInput.Range("D3").Value = 1
$wsInput.Range("D4").Value = 0.2
$wsInput.Range("D5").Value = 0.2
$wsInput.Range("D6").Value = 0.2
$wsInput.Range("D7").Value = 0.2
$wsInput.Range("D8").Value = 0.2

# Move tab selection / active cell from output to input sheet at D8
$wsInput.Select()
$wsInput.Range("D8").Select()
